$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1824561403508772
$ws.Range("C2").Value = 0.5824561403508772
$ws.Range("J2").Value = 0.01754385964912281
$ws.Range("P2").Value = 0.1228070175438596
$ws.Range("S2").Value = 0.09473684210526316
$ws.Range("C3").Value = 0.03763440860215054
$ws.Range("J3").Value = 0.03225806451612903
$ws.Range("P3").Value = 0.7741935483870968
$ws.Range("S3").Value = 0.1559139784946237
$ws.Range("J4").Value = 0.09523809523809523
$ws.Range("P4").Value = 0.7380952380952381
$ws.Range("S4").Value = 0.1666666666666667
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.06074766355140187
$ws.Range("D6").Value = 0.004672897196261682
$ws.Range("F6").Value = 0.07476635514018691
$ws.Range("J6").Value = 0.2149532710280374
$ws.Range("O6").Value = 0.01869158878504673
$ws.Range("Q6").Value = 0.1214953271028037
$ws.Range("R6").Value = 0.08878504672897196
$ws.Range("S6").Value = 0.4158878504672897
$ws.Range("B7").Value = 0.06550218340611354
$ws.Range("D7").Value = 0.01746724890829694
$ws.Range("E7").Value = 0.004366812227074236
$ws.Range("F7").Value = 0.1004366812227074
$ws.Range("J7").Value = 0.1179039301310044
$ws.Range("O7").Value = 0.008733624454148471
$ws.Range("Q7").Value = 0.1048034934497817
$ws.Range("R7").Value = 0.07860262008733625
$ws.Range("S7").Value = 0.5021834061135371
$ws.Range("B8").Value = 0.09975669099756691
$ws.Range("D8").Value = 0.0267639902676399
$ws.Range("F8").Value = 0.0583941605839416
$ws.Range("J8").Value = 0.1362530413625304
$ws.Range("O8").Value = 0.004866180048661801
$ws.Range("Q8").Value = 0.1654501216545012
$ws.Range("R8").Value = 0.0827250608272506
$ws.Range("S8").Value = 0.4257907542579075
$ws.Range("B9").Value = 0.08173076923076923
$ws.Range("D9").Value = 0.01923076923076923
$ws.Range("E9").Value = 0.004807692307692308
$ws.Range("F9").Value = 0.09134615384615384
$ws.Range("J9").Value = 0.1538461538461539
$ws.Range("O9").Value = 0.01442307692307692
$ws.Range("Q9").Value = 0.1778846153846154
$ws.Range("R9").Value = 0.08653846153846154
$ws.Range("S9").Value = 0.3701923076923077
$ws.Range("B10").Value = 0.1126651126651127
$ws.Range("D10").Value = 0.0170940170940171
$ws.Range("F10").Value = 0.05439005439005439
$ws.Range("J10").Value = 0.1460761460761461
$ws.Range("O10").Value = 0.01864801864801865
$ws.Range("Q10").Value = 0.1989121989121989
$ws.Range("R10").Value = 0.09557109557109557
$ws.Range("S10").Value = 0.3566433566433567
$ws.Range("G11").Value = 0.1761194029850746
$ws.Range("J11").Value = 0.07462686567164178
$ws.Range("K11").Value = 0.2208955223880597
$ws.Range("L11").Value = 0.5104477611940299
$ws.Range("S11").Value = 0.01791044776119403
$ws.Range("G12").Value = 0.7657142857142857
$ws.Range("J12").Value = 0.1942857142857143
$ws.Range("K12").Value = 0.01714285714285714
$ws.Range("L12").Value = 0.005714285714285714
$ws.Range("S12").Value = 0.01714285714285714
$ws.Range("G13").Value = 0.7407407407407407
$ws.Range("J13").Value = 0.2407407407407407
$ws.Range("S13").Value = 0.01851851851851852
$ws.Range("F15").Value = 0.01149425287356322
$ws.Range("H15").Value = 0.1551724137931035
$ws.Range("I15").Value = 0.07471264367816093
$ws.Range("J15").Value = 0.3448275862068966
$ws.Range("K15").Value = 0.09770114942528736
$ws.Range("M15").Value = 0.02298850574712644
$ws.Range("O15").Value = 0.04597701149425287
$ws.Range("S15").Value = 0.2471264367816092
$ws.Range("F16").Value = 0.004878048780487805
$ws.Range("H16").Value = 0.1512195121951219
$ws.Range("I16").Value = 0.07804878048780488
$ws.Range("J16").Value = 0.473170731707317
$ws.Range("K16").Value = 0.07804878048780488
$ws.Range("M16").Value = 0.02926829268292683
$ws.Range("O16").Value = 0.04878048780487805
$ws.Range("S16").Value = 0.1365853658536585
$ws.Range("F17").Value = 0.01658767772511848
$ws.Range("H17").Value = 0.1540284360189574
$ws.Range("I17").Value = 0.09715639810426541
$ws.Range("J17").Value = 0.4194312796208531
$ws.Range("K17").Value = 0.1042654028436019
$ws.Range("M17").Value = 0.01658767772511848
$ws.Range("O17").Value = 0.04502369668246445
$ws.Range("S17").Value = 0.1469194312796208
$ws.Range("F18").Value = 0.02843601895734597
$ws.Range("H18").Value = 0.1658767772511848
$ws.Range("I18").Value = 0.1042654028436019
$ws.Range("J18").Value = 0.4123222748815166
$ws.Range("K18").Value = 0.1137440758293839
$ws.Range("M18").Value = 0.02369668246445497
$ws.Range("O18").Value = 0.05213270142180094
$ws.Range("S18").Value = 0.0995260663507109
$ws.Range("F19").Value = 0.01526104417670683
$ws.Range("H19").Value = 0.2040160642570281
$ws.Range("I19").Value = 0.09397590361445783
$ws.Range("J19").Value = 0.3622489959839357
$ws.Range("K19").Value = 0.1317269076305221
$ws.Range("M19").Value = 0.02650602409638554
$ws.Range("O19").Value = 0.05783132530120482
$ws.Range("S19").Value = 0.108433734939759
